# Daily report sheet: append a new row for 2024-07-24 with zeroed totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like string into the cell as text (not an auto-converted
# date serial), then drop back to the Normal style so no extra formatting
# is left behind on the cell.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-07-24"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
